$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 9: B9 "Title, Source title, Publisher" -> "all"; C9 stays "SO, PU, SC"
$ws.Range("B9").Value = "all"
$ws.Range("C9").Value = "SO, PU, SC"

# Update the active selection shown in the saved view: was C9, now B9
$ws.Range("B9").Select()
